$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "total" row (old row 12) so the total row moves
# down to row 13. Excel carries the row-11 direct formatting (hyperlink style
# on D, shaded "Checked" style on H) onto the freshly inserted row.
$ws.Rows("12").Insert()

# Populate the new "LCD backlight" line item.
$ws.Range("C12").Value = "LCD backlight"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 100
$ws.Range("G12").Formula = "=F12*E12"

# Extend the total so it includes the new row.
$ws.Range("G13").Formula = "=SUM(G2:G12)"

# Keep the selection sensible, matching where editing left off.
$ws.Range("F13").Select()
